$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C1").EntireColumn.Delete()
Write-Host "done"
